$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    0.99999999613111856,
    0.99430578860461583,
    0.9739171963972082,
    0.96474593543879761,
    0.952506760081125,
    0.92270403050999961,
    0.91887204107486364,
    0.91306441755965384,
    0.91075849693249755,
    0.90939039230127783,
    0.90927368304865142,
    0.9093079996104958,
    0.89801900819968394,
    0.89385118492574178,
    0.8912594486716473,
    0.88875278428666671,
    0.88504466589067532,
    0.88393573032410844,
    0.99052613639868758,
    0.95933794276716866,
    0.95057254528950608,
    0.94930801857399516,
    0.96688159818500308,
    0.95385991292016503,
    0.9474027771284983,
    0.91904618196897025,
    0.91419666425008939,
    0.89271406938106135,
    0.87743899662846103,
    0.8708677837873946,
    0.86321326051440395,
    0.86153387415924443,
    0.86101385265663244
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $i + 1
    $ws.Cells.Item($row, 1).Value = $values[$i]
}
